# Insert two new data rows at the top of the "Camote" price table (rows 108-109),
# pushing all the existing rows (old 108-202) down by two (new 110-204).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(108).EntireRow.Insert()
$ws.Rows.Item(109).EntireRow.Insert()

# New row 108
$ws.Cells.Item(108, 1).Value  = 9
$ws.Cells.Item(108, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(108, 3).Value  = "Metropolitana"
$ws.Cells.Item(108, 4).Value  = 45271
$ws.Cells.Item(108, 5).Value  = 13
$ws.Cells.Item(108, 6).Value  = 100114002
$ws.Cells.Item(108, 7).Value  = "Camote"
$ws.Cells.Item(108, 8).Value  = "Sin especificar"
$ws.Cells.Item(108, 9).Value  = "Primera"
$ws.Cells.Item(108, 10).Value = 610
$ws.Cells.Item(108, 11).Value = 13000
$ws.Cells.Item(108, 12).Value = 14000
$ws.Cells.Item(108, 13).Value = 13500
$ws.Cells.Item(108, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(108, 15).Value = "Perú"
$ws.Cells.Item(108, 16).Value = 750
$ws.Cells.Item(108, 17).Value = 18
$ws.Cells.Item(108, 18).Value = "Hortaliza"

# New row 109
$ws.Cells.Item(109, 1).Value  = 9
$ws.Cells.Item(109, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(109, 3).Value  = "Metropolitana"
$ws.Cells.Item(109, 4).Value  = 45271
$ws.Cells.Item(109, 5).Value  = 13
$ws.Cells.Item(109, 6).Value  = 100114002
$ws.Cells.Item(109, 7).Value  = "Camote"
$ws.Cells.Item(109, 8).Value  = "Sin especificar"
$ws.Cells.Item(109, 9).Value  = "Primera"
$ws.Cells.Item(109, 10).Value = 970
$ws.Cells.Item(109, 11).Value = 10000
$ws.Cells.Item(109, 12).Value = 11000
$ws.Cells.Item(109, 13).Value = 10485
$ws.Cells.Item(109, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(109, 15).Value = "Perú"
$ws.Cells.Item(109, 16).Value = 582
$ws.Cells.Item(109, 17).Value = 18
$ws.Cells.Item(109, 18).Value = "Hortaliza"
